# Update the "想去人数" (column F) values across the relevant sheets
# to reflect the latest generated output (gh-pages update @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3
$ws1.Range("F3").Value = 299
$ws1.Range("F4").Value = 208
$ws1.Range("F5").Value = 2538
$ws1.Range("F6").Value = 1830
$ws1.Range("F7").Value = 351
$ws1.Range("F8").Value = 104
$ws1.Range("F9").Value = 887
$ws1.Range("F10").Value = 174

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 25

# Sheet "全部类型" (all types combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3
$ws4.Range("F3").Value = 299
$ws4.Range("F4").Value = 208
$ws4.Range("F5").Value = 2538
$ws4.Range("F6").Value = 1830
$ws4.Range("F7").Value = 351
$ws4.Range("F8").Value = 25
$ws4.Range("F9").Value = 104
$ws4.Range("F10").Value = 887
$ws4.Range("F11").Value = 174
